$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Files")

# Add new column header E1
$ws.Range("E1").Value = "*MISSING-FILE*"

# Add new rows 4-6
$ws.Range("A4").Value = "sftest.json"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1"
$ws.Range("E4").Style = "Normal"

$ws.Range("A5").Value = "test.sf.json"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1"
$ws.Range("E5").Style = "Normal"

$ws.Range("A6").Value = "test.sf.pp.json"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1"
$ws.Range("E6").Style = "Normal"
